$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Saldo" header in C1, not bold (plain style), as a new shared string
$ws.Range("C1").Value = "Saldo"
$ws.Range("C1").Font.Bold = $false

# Fill in the computed balance column
$ws.Range("C2").Value = 1200
$ws.Range("C3").Value = -200

# Update the selected cell to mimic the recorded cursor position
$ws.Range("E8").Select()
